$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'275.19"
$ws.Range("E2").Value = "'-1.43%"

$ws.Range("D3").Value = "'26.77"
$ws.Range("E3").Value = "'-2.27%"

$ws.Range("D4").Value = "'4.874"
$ws.Range("E4").Value = "'1.48%"

$ws.Range("D5").Value = "'0.06318"
$ws.Range("E5").Value = "'-0.53%"

$ws.Range("D6").Value = "'6.923"
$ws.Range("E6").Value = "'-0.37%"

$ws.Range("D7").Value = "'1.291"
$ws.Range("E7").Value = "'35.01%"

$ws.Range("D8").Value = "'0.8701"
$ws.Range("E8").Value = "'-1.01%"

$ws.Range("D9").Value = "'0.1531"
$ws.Range("E9").Value = "'3.98%"

$ws.Range("D10").Value = "'0.05058"
$ws.Range("E10").Value = "'-1.70%"

$ws.Range("D11").Value = "'0.07472"
$ws.Range("E11").Value = "'2.63%"

$ws.Range("D12").Value = "'0.02959"
$ws.Range("E12").Value = "'-5.31%"

$ws.Range("D13").Value = "'0.09059"
$ws.Range("E13").Value = "'-0.07%"

$ws.Range("D14").Value = "'0.001571"
$ws.Range("E14").Value = "'0.75%"

$ws.Range("D15").Value = "'0.0006344"
$ws.Range("E15").Value = "'0.85%"

$ws.Range("D16").Value = "'0.006034"
$ws.Range("E16").Value = "'5.22%"

$ws.Range("D17").Value = "'3.457"
$ws.Range("E17").Value = "'0.07%"

$ws.Range("D18").Value = "'3.309"
$ws.Range("E18").Value = "'-2.61%"

$ws.Range("E20").Value = "'-0.01%"

$ws.Range("D21").Value = "'0.1314"
$ws.Range("E21").Value = "'1.54%"

$ws.Range("D22").Value = "'3.926"
$ws.Range("E22").Value = "'1.35%"

$ws.Range("D23").Value = "'0.04382"
$ws.Range("E23").Value = "'1.60%"

$ws.Range("D24").Value = "'0.001170"
$ws.Range("E24").Value = "'-1.04%"

$ws.Range("E25").Value = "'-1.87%"

$ws.Range("E26").Value = "'0.08%"

$ws.Range("D27").Value = "'0.0001617"
$ws.Range("E27").Value = "'-4.40%"

$ws.Range("E40").Value = "'0.08%"

$ws.Range("D41").Value = "'0.007031"
$ws.Range("E41").Value = "'4.10%"

$ws.Range("D42").Value = "'0.1172"
$ws.Range("E42").Value = "'0.73%"

$ws.Range("D43").Value = "'0.002274"
$ws.Range("E43").Value = "'3.26%"

$ws.Range("D44").Value = "'0.01121"
$ws.Range("E44").Value = "'-10.76%"

$ws.Range("D45").Value = "'0.00005202"
$ws.Range("E45").Value = "'-0.48%"

$ws.Range("D46").Value = "'1.490"
$ws.Range("E46").Value = "'-37.38%"
